# Add a new "Portugal" market sheet, cloned from the existing "Swiss" sheet
# so it inherits the same layout/styles/merges/page setup, placed as the
# last tab (after "Swiss"), then update its market-specific content.

$wb = $excel.ActiveWorkbook

$swiss = $wb.Worksheets.Item("Swiss")
[void]$swiss.Copy($null, $swiss)

# The freshly copied sheet becomes the active sheet/tab, just like in Excel.
$portugal = $wb.ActiveSheet
$portugal.Name = "Portugal"

# Market name (B2) and Jira/ticket reference (B4) for Portugal.
$portugal.Range("B2").Value = "Portugal Market"
$portugal.Range("B4").Value = "NGC-3479/T2404"

# Match the selection left on the new sheet.
[void]$portugal.Range("B4:B5").Select()

Write-Output ("Active sheet: " + $wb.ActiveSheet.Name)
